$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing trial type labels to the "_s" (standard) variants
$ws.Range("B2:B5").Value = "cs_plus_s"
$ws.Range("B6:B9").Value = "cs_minus_s"

# Positions repeat in the same pattern as rows 2-5 / 6-9
$positions = @("(-0.5, 0.2)", "(-0.5, -0.2)", "(0.5, 0.2)", "(0.5, -0.2)")

for ($i = 0; $i -lt 4; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 1).Value = $positions[$i]
    $ws.Cells.Item($row, 2).Value = "cs_plus_ns"
}

for ($i = 0; $i -lt 4; $i++) {
    $row = 14 + $i
    $ws.Cells.Item($row, 1).Value = $positions[$i]
    $ws.Cells.Item($row, 2).Value = "cs_minus_ns"
}

# Autofit column B to match widened content
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Selection ends up on D6 after the edits
$ws.Range("D6").Select() | Out-Null
